$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 438; existing rows 438-447 shift down to 439-448.
$ws.Rows(438).Insert()

# Populate the newly inserted row 438 with the new price record.
$ws.Range("A438").Value = 10
$ws.Range("B438").Value = "Vega Modelo de Temuco"
$ws.Range("C438").Value = "La Araucanía"
$ws.Range("D438").Value = 45239
$ws.Range("E438").Value = 9
$ws.Range("F438").Value = "Fruta"
$ws.Range("G438").Value = 100103
$ws.Range("H438").Value = "Frutos de hueso (carozo)"
$ws.Range("I438").Value = 100103004
$ws.Range("J438").Value = "Durazno"
$ws.Range("K438").Value = "Florida King"
$ws.Range("L438").Value = "Primera"
$ws.Range("M438").Value = 185
$ws.Range("N438").Value = 20000
$ws.Range("O438").Value = 20000
$ws.Range("P438").Value = 20000
$ws.Range("Q438").Value = "$/bandeja 10 kilos granel"
$ws.Range("R438").Value = "Provincia de Limarí"
$ws.Range("S438").Value = 2000
$ws.Range("T438").Value = 10
